$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Some Price values (e.g. "0.520", "21.70") look numeric but must stay as
# plain text so trailing zeros / exact digits are preserved (matches source feed).
# We flip NumberFormat to Text ("@") before assigning such values, then restore
# the cell style back to Normal/General so no stray formatting is left behind.

$ws.Range('D2').Value = '27.216.14'
$ws.Range('E2').Value = '  +0.89%  '
$ws.Range('D3').Value = '1.684.96'
$ws.Range('E3').Value = '  +0.48%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '215.79'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.27%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.520'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +0.31%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('E8').Value = '  +2.17%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '21.70'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +6.79%  '
$ws.Range('E10').Value = '  +0.69%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0889'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('D12').Value = '1.921.78'
$ws.Range('E12').Value = '  +0.43%  '
$ws.Range('D13').Value = '1.675.26'
$ws.Range('E13').Value = '  -1.28%  '
$ws.Range('E14').Value = '  +1.67%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.545'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +2.98%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '66.43'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.84%  '
$ws.Range('D17').Value = '27.211.28'
$ws.Range('E17').Value = '  +0.81%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '239.44'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.94%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '8.10'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.27%  '
$ws.Range('D20').Value = '0.0₃0743'
$ws.Range('E20').Value = '  +1.43%  '
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('E22').Value = '  +2.66%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '9.53'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +3.59%  '
$ws.Range('E24').Value = '  -3.45%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '148.41'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +1.92%  '
$ws.Range('E26').Value = '  +0.57%  '
$ws.Range('E27').Value = '  +1.70%  '
$ws.Range('E28').Value = '  +0.86%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.999'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('E30').Value = '  +0.20%  '
$ws.Range('E31').Value = '  +0.58%  '
$ws.Range('D32').Value = '1.579.25'
$ws.Range('E32').Value = '  +6.24%  '
$ws.Range('E33').Value = '  +1.72%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '3.24'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +2.64%  '
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.603'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +3.06%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.945'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +4.67%  '
$ws.Range('E38').Value = '  -0.85%  '
$ws.Range('E39').Value = '  -0.46%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.07'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +4.25%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '69.26'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +2.62%  '
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('E43').Value = '  -4.22%  '
$ws.Range('E44').Value = '  -2.63%  '
$ws.Range('D45').Value = '1.830.57'
$ws.Range('E45').Value = '  +0.64%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.788'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +0.83%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '90.95'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.35%  '
$ws.Range('E48').Value = '  +3.93%  '
$ws.Range('E49').Value = '  +1.31%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '8.17'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +6.17%  '
$ws.Range('E51').Value = '  +2.04%  '
